# Automatische test-sync: 2025-06-19 14:30:10
# Adds the new incoming-mail log row (row 14) to the "Logs" sheet, extends
# the conditional-formatting ranges to cover it, and bumps the "Overig"
# tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row ---------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A14").Value = "Vragen over samenwerking"
$logs.Range("B14").Value = "mailmind.test@zohomail.eu"
$logs.Range("C14").Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D14").Value = "Overig"
$logs.Range("F14").Value = "2025-06-19 14:28:10"
$logs.Range("G14").Value = "Nee"

# --- Extend the conditional formatting ranges to include the new row ------
$catRules = $logs.Range("D2:D13").FormatConditions
$catCount = $catRules.Count()
for ($i = 1; $i -le $catCount; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D14"))
}

$answeredRules = $logs.Range("G2:G13").FormatConditions
$answeredCount = $answeredRules.Count()
for ($i = 1; $i -le $answeredCount; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G14"))
}

# --- Dashboard sheet: bump the "Overig" count from 6 to 7 ------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 7
